# netCrypto.xlsx daily update: "Add files via upload"
#
# The underlying change is a single data edit on sheet "SheetName1":
# cell T2 (USD Amount for the Roobic/Crypto deposit row) is updated
# from 93098 to 90764.
#
# (The diff also shows the workbook's x15ac:absPath save-folder date,
# the xr:revisionPtr session documentId GUID, and the bookViews
# xWindow screen-position — these are environment/session artifacts
# that Excel stamps on its own from the machine's file path, a fresh
# session GUID, and the OS window position at save time; they carry
# no workbook data and are not exposed through the Excel object model
# for scripts to set, so they are left to the host application.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

$ws.Range("T2").Value = 90764
